# Commit: "Submitting changes applied to backstage scenarios."
#
# Cell A5 contains a wrapped, multi-line verification-type list. A new
# line "USV+" is inserted right before "Hearing", which grows the cell
# to one more wrapped line, so the row height grows accordingly
# (116 -> 130.5, i.e. +1 default row height of 14.5).
#
# Cell A15 ("Appeal Decision") is left content-wise unchanged.
#
# The active selection also moved from A15 to A4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newText = "                                " + "`n" + `
    "                                ARSV" + "`n" + `
    "                                ASV" + "`n" + `
    "                                USV" + "`n" + `
    "                                VSV" + "`n" + `
    "                                NSV" + "`n" + `
    "                                USV+" + "`n" + `
    "                                Hearing" + "`n" + `
    "                            "

$ws.Range("A5").Value = $newText
$ws.Rows.Item(5).RowHeight = 130.5

$ws.Range("A4").Select()
